$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert rows from bottom to top so earlier row references remain valid
$ws.Rows("21:21").Insert()
$ws.Rows("18:18").Insert()
$ws.Rows("15:15").Insert()
$ws.Rows("12:12").Insert()
$ws.Rows("9:9").Insert()
$ws.Rows("6:6").Insert()

# Fill in new o3-mini rows for each dataset group
$ws.Range("B6").Value = "o3-mini"
$ws.Range("C6").Value = 154
$ws.Range("D6").Value = 0.59
$ws.Range("E6").Value = 0.69

$ws.Range("B10").Value = "o3-mini"
$ws.Range("C10").Value = 150
$ws.Range("D10").Value = 0.57
$ws.Range("E10").Value = 0.74

$ws.Range("B14").Value = "o3-mini"
$ws.Range("C14").Value = 213
$ws.Range("D14").Value = 0.51
$ws.Range("E14").Value = 0.69

$ws.Range("B18").Value = "o3-mini"
$ws.Range("C18").Value = 143
$ws.Range("D18").Value = 0.53
$ws.Range("E18").Value = 0.78

$ws.Range("B22").Value = "o3-mini"
$ws.Range("C22").Value = 178
$ws.Range("D22").Value = 0.48
$ws.Range("E22").Value = 0.7

$ws.Range("B26").Value = "o3-mini"
$ws.Range("C26").Value = 140
$ws.Range("D26").Value = 0.5
$ws.Range("E26").Value = 0.75

# New text annotations in column G
$ws.Range("G19").Value = "10.35% yes"
$ws.Range("G19").NumberFormat = "0.00%"
$ws.Range("G21").Value = "5.77% yes"

# Update selection to match the final state
$ws.Range("G10").Select()

Write-Host "done"
